# Add a "2022-Q3" sheet (fund holdings for the new quarter) right after the
# "总计" (summary) sheet, and add a corresponding summary row on "总计".
#
# Helper: write a text value into a cell while forcing the "text" cell type
# (Excel otherwise auto-coerces numeric-looking strings like "7.18" or
# "000339" into numbers, which would lose leading zeros / change cell type).
# NumberFormat="@" forces text-on-entry; resetting Style back to "Normal"
# afterwards avoids leaving a stray custom cell style behind.
function Set-TextValue($ws, $row, $col, $val) {
    $cell = $ws.Cells.Item($row, $col)
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = "Normal"
}

function Set-FundRow($ws, $r, $code, $name, $scale, $totalPos, $posPct, $mktValue, $rank) {
    Set-TextValue $ws $r 2 $code
    Set-TextValue $ws $r 3 $name
    Set-TextValue $ws $r 4 $scale
    Set-TextValue $ws $r 5 $totalPos
    Set-TextValue $ws $r 6 $posPct
    Set-TextValue $ws $r 7 $mktValue
    $ws.Cells.Item($r, 8).Value = $rank
}

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Create the new "2022-Q3" worksheet right after "总计".
# ---------------------------------------------------------------------
$sheetTotal = $wb.Worksheets.Item(1)
$newSheet = $wb.Worksheets.Add($null, $sheetTotal)
$newSheet.Name = "2022-Q3"

# Use the (still unshifted at this point) "2022-Q2" sheet as a formatting
# template, since its layout (header row style + data row/column styles) is
# identical to what the new quarter sheet needs.
$template = $wb.Worksheets.Item("2022-Q2")

# Header row formatting (bold / centered / bordered style).
$template.Range("B1:H1").Copy()
$newSheet.Range("B1:H1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Data-row formatting (plain style) across all 13 data rows.
$template.Range("B2:H2").Copy()
$newSheet.Range("B2:H14").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Index-column formatting (bold/centered "A" column) across all 13 rows.
$template.Range("A2").Copy()
$newSheet.Range("A2:A14").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Header values.
$newSheet.Range("B1").Value = "基金代码"
$newSheet.Range("C1").Value = "基金名称"
$newSheet.Range("D1").Value = "基金规模"
$newSheet.Range("E1").Value = "股票总仓位"
$newSheet.Range("F1").Value = "仓位占比"
$newSheet.Range("G1").Value = "持有市值(亿元)"
$newSheet.Range("H1").Value = "仓位排名"

# Index column (0-based row counter).
for ($i = 0; $i -le 12; $i++) {
    $newSheet.Cells.Item($i + 2, 1).Value = $i
}

# Fund holdings data rows.
Set-FundRow $newSheet 2  "000339" "长城医疗保健混合A" "7.18" "87.41" "2.98" "0.2140" 8
Set-FundRow $newSheet 3  "011673" "长城医药科技六个月持有期混合型证券投资基金A" "5.75" "87.90" "3.02" "0.1736" 8
Set-FundRow $newSheet 4  "005689" "中银医疗保健灵活配置混合A" "5.59" "89.04" "3.04" "0.1699" 9
Set-FundRow $newSheet 5  "008786" "长城健康生活灵活配置混合" "6.21" "79.23" "2.51" "0.1559" 6
Set-FundRow $newSheet 6  "013293" "长城健康消费混合" "5.90" "76.79" "2.02" "0.1192" 10
Set-FundRow $newSheet 7  "013037" "长城大健康混合A" "5.40" "72.73" "2.18" "0.1177" 9
Set-FundRow $newSheet 8  "007718" "中银创新医疗混合A" "2.82" "92.35" "3.12" "0.0880" 10
Set-FundRow $newSheet 9  "013441" "西藏东财创新医疗六个月定开混合" "0.49" "82.53" "4.71" "0.0231" 8
Set-FundRow $newSheet 10 "011674" "长城医药科技六个月持有期混合型证券投资基金C" "0.54" "87.90" "3.02" "0.0163" 8
Set-FundRow $newSheet 11 "010159" "中银医疗保健灵活配置混合C" "0.32" "89.04" "3.04" "0.0097" 9
Set-FundRow $newSheet 12 "013038" "长城大健康混合C" "0.26" "72.73" "2.18" "0.0057" 9
Set-FundRow $newSheet 13 "010500" "中银创新医疗混合C" "0.18" "92.35" "3.12" "0.0056" 10
Set-FundRow $newSheet 14 "015562" "长城医疗保健混合C" "0.13" "87.41" "2.98" "0.0039" 8

# ---------------------------------------------------------------------
# 2. Insert a new summary row for "2022-Q3" at the top of the "总计" data
#    (row 2), pushing the existing quarters down by one row.
# ---------------------------------------------------------------------
$totalSheet = $wb.Worksheets.Item("总计")
$totalSheet.Rows(2).Insert()

# The Insert() above leaves A2 without its usual style and B2:D2 with a
# style inherited from the header row; restore the plain styling used by
# the rest of the table by copying formats from row 3 (the old row 2,
# now shifted down, which still has the correct original styling).
$totalSheet.Range("A3").Copy()
$totalSheet.Range("A2").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$totalSheet.Range("B3:D3").Copy()
$totalSheet.Range("B2:D2").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$totalSheet.Cells.Item(2, 1).Value = 0
$totalSheet.Cells.Item(2, 2).Value = "2022-Q3"
$totalSheet.Cells.Item(2, 3).Value = 13
$totalSheet.Cells.Item(2, 4).Value = 1.1

# Renumber the 0-based index column for the rows that shifted down.
for ($r = 3; $r -le 9; $r++) {
    $totalSheet.Cells.Item($r, 1).Value = $r - 2
}

Write-Host "2022-Q3 sheet added and 总计 summary updated"
